$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "IsHuman"
$ws.Range("E1").Value = "IsEsh"

$ws.Range("D2").Value = $true
$ws.Range("E2").Value = 0

$ws.Range("E3").Select()
